# Add a new "2022-Q1" sheet (with fund holding detail) before the "总计"
# (summary) sheet, and add a corresponding "2022-Q1" row at the top of the
# "总计" summary sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: style a header-row cell (bold, bordered, centered) to visually
# match the other header rows already present in the workbook.
# ---------------------------------------------------------------------
function Style-Header($cell) {
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Helper: style the small numeric index cells in column A (bold, bordered,
# centered) the same way the existing sheets do.
function Style-IndexCell($cell) {
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Helper: force a cell to be written as literal text (not auto-converted
# to a number), then clear the number-format override so no stray style
# is left behind on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Locate the existing "总计" sheet and insert the new "2022-Q1" sheet
#    immediately before it.
#
#    NOTE: after Worksheets.Add(...) / renaming, the original variable
#    can end up pointing at the wrong worksheet (the runtime appears to
#    resolve worksheet variables positionally), so the "总计" sheet is
#    re-fetched by name below, once the new sheet has been created, to
#    make sure later writes land on the correct sheet.
# ---------------------------------------------------------------------
$zongjiSheetBefore = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Add($zongjiSheetBefore)
$q1Sheet.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 2. Fill in the "2022-Q1" sheet headers and data.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q1Sheet.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    Style-Header $cell
}

$rows = @(
    @("001933", "华商新兴活力灵活配置混合", "6.91", "92.28", "4.02", "0.2778", 9),
    @("630008", "华商策略精选混合",         "4.71", "76.99", "3.66", "0.1724", 3),
    @("011851", "天弘先进制造混合型证券投资基金A", "2.72", "91.41", "5.03", "0.1368", 9),
    @("011852", "天弘先进制造混合型证券投资基金C", "0.70", "91.41", "5.03", "0.0352", 9)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    $idxCell = $q1Sheet.Cells.Item($rowNum, 1)
    $idxCell.Value = $r
    Style-IndexCell $idxCell

    Set-TextValue ($q1Sheet.Cells.Item($rowNum, 2)) $data[0]
    Set-TextValue ($q1Sheet.Cells.Item($rowNum, 3)) $data[1]
    Set-TextValue ($q1Sheet.Cells.Item($rowNum, 4)) $data[2]
    Set-TextValue ($q1Sheet.Cells.Item($rowNum, 5)) $data[3]
    Set-TextValue ($q1Sheet.Cells.Item($rowNum, 6)) $data[4]
    Set-TextValue ($q1Sheet.Cells.Item($rowNum, 7)) $data[5]

    $q1Sheet.Cells.Item($rowNum, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 3. Insert a new summary row at the top of the "总计" sheet's data table
#    (row 2) for "2022-Q1", pushing the existing rows down.
# ---------------------------------------------------------------------
$zongjiSheet = $wb.Worksheets.Item("总计")
$zongjiSheet.Rows.Item(2).Insert()

$idxCell = $zongjiSheet.Cells.Item(2, 1)
$idxCell.Value = 0
Style-IndexCell $idxCell

$dateCell = $zongjiSheet.Cells.Item(2, 2)
$dateCell.Value = "2022-Q1"
$dateCell.Style = "Normal"

$zongjiSheet.Cells.Item(2, 3).Value = 4
$zongjiSheet.Cells.Item(2, 3).Style = "Normal"
$zongjiSheet.Cells.Item(2, 4).Value = 0.62
$zongjiSheet.Cells.Item(2, 4).Style = "Normal"

# The row-insert operation pushes the previously existing rows down but
# keeps their original index values (0, 1, 2) in column A. Renumber them
# sequentially (1, 2, 3) to account for the newly inserted row above.
$zongjiSheet.Cells.Item(3, 1).Value = 1
$zongjiSheet.Cells.Item(4, 1).Value = 2
$zongjiSheet.Cells.Item(5, 1).Value = 3

Write-Host "2022-Q1 sheet and summary row added."
